# Mock-Up.pptx tweaks:
#  1. The "auto-updating" date placeholder shown on the slide master and every
#     slide layout reads 27.02.2022 -> bump it to 28.02.2022.
#  2. On the last slide (the "Stundensatz - Berechnung" mock-up), the label
#     above the 100%-workload gross-wage figure (shape "Textfeld 55") changes
#     from "Bruttolohn Pensum 100%:" to "Berechneter Stundensatz:". (There is
#     a second, unrelated shape further down with the same original caption
#     that must stay untouched.)

$p = $ppt.ActivePresentation

$oldDate = "27.02.2022"
$newDate = "28.02.2022"

# --- Slide master: update the date field text ---
$master = $p.SlideMaster
foreach ($shape in $master.Shapes) {
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every slide layout off that master: same date field ---
foreach ($layout in $master.CustomLayouts) {
    foreach ($shape in $layout.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Last slide: rename the specific "Bruttolohn Pensum 100%:" label ---
$lastSlide = $p.Slides.Item($p.Slides.Count)
$label = $lastSlide.Shapes.Item("Textfeld 55")
if ($label.TextFrame.TextRange.Text -eq "Bruttolohn Pensum 100%:") {
    $label.TextFrame.TextRange.Text = "Berechneter Stundensatz:"
}
